$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (copy the header style from H1 so I1/J1 match the
# bold/centered/bordered header formatting used by the rest of row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for rows 2-4
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 2
